# ============================================================================
# Edit: add "Player Info" sheet (before "ODI Batting"), add "ODI Batting Extra"
# sheet (after "ODI Bowling"), and on the existing "ODI Batting"/"ODI Bowling"
# sheets rename MATCH_CARD_LINK -> MATCH_CODE, replacing the full scorecard
# URL with just the numeric match code extracted from it. Also drop a handful
# of stray empty INNING_NUMBER cells on "ODI Batting".
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1. Insert "Player Info" before "ODI Batting"
# ----------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$playerInfo.Range("A2").Value = "'4222"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Michael A Leask"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# ----------------------------------------------------------------------
# 2. Append "ODI Batting Extra" after "ODI Bowling" (last sheet)
# ----------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $extraHeaders.Length; $col++) {
    $battingExtra.Cells.Item(1, $col).Value = $extraHeaders[$col - 1]
}
$extraHeaderRange = $battingExtra.Range("A1:F1")
$extraHeaderRange.Font.Bold = $true
$extraHeaderRange.Borders.LineStyle = 1
$extraHeaderRange.HorizontalAlignment = -4108
$extraHeaderRange.VerticalAlignment = -4160

# MATCH_CODE(A, text), BATTING_POSITION(B, number|blank), NUM_4(C, text),
# NUM_6(D, text), PERCENT_RUNS_OF_TOTAL(E, text), MAN_OF_MATCH(F, text)
$extraRows = @(
    @("4575", 7, "1", "0", "2.91%", "NO"),
    @("4576", 7, "1", "0", "2.67%", "NO"),
    @("4578", 7, "4", "1", "21.66%", "NO"),
    @("4581", 7, "0", "0", "2.92%", "NO"),
    @("4604", 7, "1", "2", "6.98%", "NO"),
    @("4610", 6, "1", "0", "2.78%", "NO"),
    @("4612", 7, "0", "0", "", "NO"),
    @("4625", 7, "9", "4", "27.78%", "NO"),
    @("4629", 7, "4", "0", "14.89%", "NO"),
    @("4631", $null, "", "", "", "NO"),
    @("4632", 7, "0", "0", "1.97%", "NO"),
    @("4635", $null, "", "", "", "NO"),
    @("4677", 7, "2", "1", "8.45%", "NO"),
    @("4681", 7, "3", "0", "15.38%", "NO"),
    @("4680", 8, "0", "1", "6.73%", "NO"),
    @("4684", 8, "", "", "", "NO"),
    @("4702", 5, "", "", "", "NO"),
    @("4703", 7, "3", "9", "39.05%", "NO"),
    @("4705", 8, "1", "2", "13.12%", "NO"),
    @("4706", $null, "", "", "", "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $battingExtra.Cells.Item($r, 1).Value = "'" + $row[0]
    $battingExtra.Cells.Item($r, 1).Style = "Normal"

    if ($null -ne $row[1]) {
        $battingExtra.Cells.Item($r, 2).Value = $row[1]
    } else {
        $battingExtra.Cells.Item($r, 2).Value = "'"
        $battingExtra.Cells.Item($r, 2).Style = "Normal"
    }

    $battingExtra.Cells.Item($r, 3).Value = "'" + $row[2]
    $battingExtra.Cells.Item($r, 3).Style = "Normal"

    $battingExtra.Cells.Item($r, 4).Value = "'" + $row[3]
    $battingExtra.Cells.Item($r, 4).Style = "Normal"

    $battingExtra.Cells.Item($r, 5).Value = "'" + $row[4]
    $battingExtra.Cells.Item($r, 5).Style = "Normal"

    $battingExtra.Cells.Item($r, 6).Value = $row[5]

    $r = $r + 1
}

# ----------------------------------------------------------------------
# 3. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, strip url down to code,
#    drop a few stray empty INNING_NUMBER (column B) cells.
#
# NOTE: re-fetch sheet objects by name here rather than reusing the
# references captured at the top of the script -- `Worksheets.Add(Before)`
# re-purposes the identity bound to the `Before` argument's variable to
# point at the freshly inserted sheet, so `$battingSheet`/`$bowlingSheet`
# from before the Add calls above no longer point at the intended sheets.
# ----------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$battingSheet.Range("D1").Value = "MATCH_CODE"

for ($row = 2; $row -le 63; $row++) {
    $cell = $battingSheet.Cells.Item($row, 4)
    $text = $cell.Text
    if ($text -match "MatchCode=(\d+)") {
        $code = $matches[1]
        $cell.Value = "'" + $code
        $cell.Style = "Normal"
    }
}

$emptyInningRows = @(7, 12, 15, 22, 23, 26, 36, 39, 40, 59, 60)
foreach ($row in $emptyInningRows) {
    $battingSheet.Cells.Item($row, 2).ClearContents()
}

# ----------------------------------------------------------------------
# 4. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, strip url down to code.
# ----------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

for ($row = 2; $row -le 57; $row++) {
    $cell = $bowlingSheet.Cells.Item($row, 2)
    $text = $cell.Text
    if ($text -match "MatchCode=(\d+)") {
        $code = $matches[1]
        $cell.Value = "'" + $code
        $cell.Style = "Normal"
    }
}
